$d = $word.ActiveDocument

# [Music] -> [Muziki]
$d.Content.Find.Execute("[Music]", $true, $true, $false, $false, $false,
                         $true, 1, $false, "[Muziki]", 2)

# for example -> kwa mfano (all occurrences)
$d.Content.Find.Execute("for example", $true, $true, $false, $false, $false,
                         $true, 1, $false, "kwa mfano", 2)

# [PAUSE] -> [SItisha] (all occurrences)
$d.Content.Find.Execute("[PAUSE]", $true, $true, $false, $false, $false,
                         $true, 1, $false, "[SItisha]", 2)
